$wb = $excel.ActiveWorkbook
$first = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($first)
$ws.Name = "E"

# Write column A (instance names) first so the new shared-strings are interned
# in the same order as the target workbook: steinb01..steinb18 before best/l/r.
$ws.Cells.Item(2, 1).Value = "\\steinb01"
$ws.Cells.Item(3, 1).Value = "\\steinb02"
$ws.Cells.Item(4, 1).Value = "\\steinb03"
$ws.Cells.Item(5, 1).Value = "\\steinb04"
$ws.Cells.Item(6, 1).Value = "\\steinb05"
$ws.Cells.Item(7, 1).Value = "\\steinb06"
$ws.Cells.Item(8, 1).Value = "\\steinb07"
$ws.Cells.Item(9, 1).Value = "\\steinb08"
$ws.Cells.Item(10, 1).Value = "\\steinb09"
$ws.Cells.Item(11, 1).Value = "\\steinb10"
$ws.Cells.Item(12, 1).Value = "\\steinb11"
$ws.Cells.Item(13, 1).Value = "\\steinb12"
$ws.Cells.Item(14, 1).Value = "\\steinb13"
$ws.Cells.Item(15, 1).Value = "\\steinb14"
$ws.Cells.Item(16, 1).Value = "\\steinb15"
$ws.Cells.Item(17, 1).Value = "\\steinb16"
$ws.Cells.Item(18, 1).Value = "\\steinb17"
$ws.Cells.Item(19, 1).Value = "\\steinb18"

# Header row
$ws.Cells.Item(1, 1).Value = "instance"
$ws.Cells.Item(1, 2).Value = "min"
$ws.Cells.Item(1, 3).Value = "max"
$ws.Cells.Item(1, 4).Value = "avg"
$ws.Cells.Item(1, 5).Value = "median"
$ws.Cells.Item(1, 6).Value = "time"
$ws.Cells.Item(1, 7).Value = "best"
$ws.Cells.Item(1, 8).Value = "l"
$ws.Cells.Item(1, 9).Value = "l"
$ws.Cells.Item(1, 10).Value = "r"
$ws.Cells.Item(1, 11).Value = "r"
$ws.Cells.Item(1, 12).Value = "solved"

# Data rows 2-19 (columns B-L)
$ws.Cells.Item(2, 2).Value = 83
$ws.Cells.Item(2, 3).Value = 83
$ws.Cells.Item(2, 4).Value = 83
$ws.Cells.Item(2, 5).Value = 83
$ws.Cells.Item(2, 6).Value = 163.8
$ws.Cells.Item(2, 7).Value = 91.7
$ws.Cells.Item(2, 8).Value = 7
$ws.Cells.Item(2, 9).Value = 7
$ws.Cells.Item(2, 10).Value = 2
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 10
$ws.Cells.Item(3, 2).Value = 83
$ws.Cells.Item(3, 3).Value = 83
$ws.Cells.Item(3, 4).Value = 83
$ws.Cells.Item(3, 5).Value = 83
$ws.Cells.Item(3, 6).Value = 124.5
$ws.Cells.Item(3, 7).Value = 14.6
$ws.Cells.Item(3, 8).Value = 7
$ws.Cells.Item(3, 9).Value = 7
$ws.Cells.Item(3, 10).Value = 3
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 10
$ws.Cells.Item(4, 2).Value = 138
$ws.Cells.Item(4, 3).Value = 138
$ws.Cells.Item(4, 4).Value = 138
$ws.Cells.Item(4, 5).Value = 138
$ws.Cells.Item(4, 6).Value = 218.7
$ws.Cells.Item(4, 7).Value = 95.9
$ws.Cells.Item(4, 8).Value = 8
$ws.Cells.Item(4, 9).Value = 8
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 10
$ws.Cells.Item(5, 2).Value = 59
$ws.Cells.Item(5, 3).Value = 59
$ws.Cells.Item(5, 4).Value = 59
$ws.Cells.Item(5, 5).Value = 59
$ws.Cells.Item(5, 6).Value = 129.3
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 7
$ws.Cells.Item(5, 9).Value = 7
$ws.Cells.Item(5, 10).Value = 3
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 10
$ws.Cells.Item(6, 2).Value = 61
$ws.Cells.Item(6, 3).Value = 61
$ws.Cells.Item(6, 4).Value = 61
$ws.Cells.Item(6, 5).Value = 61
$ws.Cells.Item(6, 6).Value = 253.5
$ws.Cells.Item(6, 7).Value = 128.2
$ws.Cells.Item(6, 8).Value = 5
$ws.Cells.Item(6, 9).Value = 6
$ws.Cells.Item(6, 10).Value = 4
$ws.Cells.Item(6, 11).Value = 4
$ws.Cells.Item(6, 12).Value = 10
$ws.Cells.Item(7, 2).Value = 125
$ws.Cells.Item(7, 3).Value = 125
$ws.Cells.Item(7, 4).Value = 125
$ws.Cells.Item(7, 5).Value = 125
$ws.Cells.Item(7, 6).Value = 330.5
$ws.Cells.Item(7, 7).Value = 191.1
$ws.Cells.Item(7, 8).Value = 7
$ws.Cells.Item(7, 9).Value = 7
$ws.Cells.Item(7, 10).Value = 5
$ws.Cells.Item(7, 11).Value = 5
$ws.Cells.Item(7, 12).Value = 10
$ws.Cells.Item(8, 2).Value = 111
$ws.Cells.Item(8, 3).Value = 112
$ws.Cells.Item(8, 4).Value = 111.1
$ws.Cells.Item(8, 5).Value = 111
$ws.Cells.Item(8, 6).Value = 237.5
$ws.Cells.Item(8, 7).Value = 100.6
$ws.Cells.Item(8, 8).Value = 7
$ws.Cells.Item(8, 9).Value = 7
$ws.Cells.Item(8, 10).Value = 5
$ws.Cells.Item(8, 11).Value = 6
$ws.Cells.Item(8, 12).Value = 10
$ws.Cells.Item(9, 2).Value = 107
$ws.Cells.Item(9, 3).Value = 107
$ws.Cells.Item(9, 4).Value = 107
$ws.Cells.Item(9, 5).Value = 107
$ws.Cells.Item(9, 6).Value = 140.9
$ws.Cells.Item(9, 7).Value = 14.1
$ws.Cells.Item(9, 8).Value = 8
$ws.Cells.Item(9, 9).Value = 8
$ws.Cells.Item(9, 10).Value = 4
$ws.Cells.Item(9, 11).Value = 4
$ws.Cells.Item(9, 12).Value = 10
$ws.Cells.Item(10, 2).Value = 221
$ws.Cells.Item(10, 3).Value = 221
$ws.Cells.Item(10, 4).Value = 221
$ws.Cells.Item(10, 5).Value = 221
$ws.Cells.Item(10, 6).Value = 165.6
$ws.Cells.Item(10, 7).Value = 0.5
$ws.Cells.Item(10, 8).Value = 9
$ws.Cells.Item(10, 9).Value = 10
$ws.Cells.Item(10, 10).Value = 4
$ws.Cells.Item(10, 11).Value = 4
$ws.Cells.Item(10, 12).Value = 10
$ws.Cells.Item(11, 2).Value = 88
$ws.Cells.Item(11, 3).Value = 88
$ws.Cells.Item(11, 4).Value = 88
$ws.Cells.Item(11, 5).Value = 88
$ws.Cells.Item(11, 6).Value = 342.2
$ws.Cells.Item(11, 7).Value = 159.6
$ws.Cells.Item(11, 8).Value = 5
$ws.Cells.Item(11, 9).Value = 5
$ws.Cells.Item(11, 10).Value = 7
$ws.Cells.Item(11, 11).Value = 7
$ws.Cells.Item(11, 12).Value = 10
$ws.Cells.Item(12, 2).Value = 88
$ws.Cells.Item(12, 3).Value = 88
$ws.Cells.Item(12, 4).Value = 88
$ws.Cells.Item(12, 5).Value = 88
$ws.Cells.Item(12, 6).Value = 380.7
$ws.Cells.Item(12, 7).Value = 194.9
$ws.Cells.Item(12, 8).Value = 3
$ws.Cells.Item(12, 9).Value = 4
$ws.Cells.Item(12, 10).Value = 4
$ws.Cells.Item(12, 11).Value = 4
$ws.Cells.Item(12, 12).Value = 10
$ws.Cells.Item(13, 2).Value = 174
$ws.Cells.Item(13, 3).Value = 176
$ws.Cells.Item(13, 4).Value = 174.9
$ws.Cells.Item(13, 5).Value = 175
$ws.Cells.Item(13, 6).Value = 513
$ws.Cells.Item(13, 7).Value = 327.1
$ws.Cells.Item(13, 8).Value = 9
$ws.Cells.Item(13, 9).Value = 11
$ws.Cells.Item(13, 10).Value = 2
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 10
$ws.Cells.Item(14, 2).Value = 173
$ws.Cells.Item(14, 3).Value = 174
$ws.Cells.Item(14, 4).Value = 173.9
$ws.Cells.Item(14, 5).Value = 174
$ws.Cells.Item(14, 6).Value = 261.5
$ws.Cells.Item(14, 7).Value = 71.4
$ws.Cells.Item(14, 8).Value = 14
$ws.Cells.Item(14, 9).Value = 17
$ws.Cells.Item(14, 10).Value = 5
$ws.Cells.Item(14, 11).Value = 5
$ws.Cells.Item(14, 12).Value = 10
$ws.Cells.Item(15, 2).Value = 239
$ws.Cells.Item(15, 3).Value = 239
$ws.Cells.Item(15, 4).Value = 239
$ws.Cells.Item(15, 5).Value = 239
$ws.Cells.Item(15, 6).Value = 490.7
$ws.Cells.Item(15, 7).Value = 261.8
$ws.Cells.Item(15, 8).Value = 13
$ws.Cells.Item(15, 9).Value = 15
$ws.Cells.Item(15, 10).Value = 11
$ws.Cells.Item(15, 11).Value = 11
$ws.Cells.Item(15, 12).Value = 10
$ws.Cells.Item(16, 2).Value = 325
$ws.Cells.Item(16, 3).Value = 328
$ws.Cells.Item(16, 4).Value = 325.75
$ws.Cells.Item(16, 5).Value = 325
$ws.Cells.Item(16, 6).Value = 688.5
$ws.Cells.Item(16, 7).Value = 417.75
$ws.Cells.Item(16, 8).Value = 12
$ws.Cells.Item(16, 9).Value = 14
$ws.Cells.Item(16, 10).Value = 6
$ws.Cells.Item(16, 11).Value = 6
$ws.Cells.Item(16, 12).Value = 4
$ws.Cells.Item(17, 2).Value = 127
$ws.Cells.Item(17, 3).Value = 130
$ws.Cells.Item(17, 4).Value = 128.5
$ws.Cells.Item(17, 5).Value = 128.5
$ws.Cells.Item(17, 6).Value = 518.1
$ws.Cells.Item(17, 7).Value = 279.2
$ws.Cells.Item(17, 8).Value = 10
$ws.Cells.Item(17, 9).Value = 16
$ws.Cells.Item(17, 10).Value = 2
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 10
$ws.Cells.Item(18, 2).Value = 132
$ws.Cells.Item(18, 3).Value = 142
$ws.Cells.Item(18, 4).Value = 133.4
$ws.Cells.Item(18, 5).Value = 132
$ws.Cells.Item(18, 6).Value = 545.4
$ws.Cells.Item(18, 7).Value = 418.4
$ws.Cells.Item(18, 8).Value = 8
$ws.Cells.Item(18, 9).Value = 11
$ws.Cells.Item(18, 10).Value = 6
$ws.Cells.Item(18, 11).Value = 8
$ws.Cells.Item(18, 12).Value = 10
$ws.Cells.Item(19, 2).Value = 218
$ws.Cells.Item(19, 3).Value = 219
$ws.Cells.Item(19, 4).Value = 218.9
$ws.Cells.Item(19, 5).Value = 219
$ws.Cells.Item(19, 6).Value = 514
$ws.Cells.Item(19, 7).Value = 268.9
$ws.Cells.Item(19, 8).Value = 10
$ws.Cells.Item(19, 9).Value = 10
$ws.Cells.Item(19, 10).Value = 8
$ws.Cells.Item(19, 11).Value = 8
$ws.Cells.Item(19, 12).Value = 10

# Row 20 average formulas
$ws.Range("F20").Formula = "=AVERAGE(F2:F19)"
$ws.Range("G20").Formula = "=AVERAGE(G2:G19)"

# Scientific-notation number format on specific highlighted cells (reuses existing style)
$ws.Range("G4").NumberFormat = "0.00E+00"
$ws.Range("H7").NumberFormat = "0.00E+00"
$ws.Range("G9").NumberFormat = "0.00E+00"

# Column A width
$ws.Columns.Item(1).ColumnWidth = 16.28515625

# View: select G14 and make this new sheet the active/selected tab
$ws.Range("G14").Select()
$ws.Select()